$d = $word.ActiveDocument

# 1. Update the letter date from September 19, 2025 to September 21, 2025
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing-address line "195 Lewis Road Suite, San Jose CA 95111"
#    into two separate paragraphs:
#      "195 Lewis Road Suite"
#      "San Jose, CA 95111"
$addr = $d.Content
$found = $addr.Find.Execute("195 Lewis Road Suite, San Jose CA 95111", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $addr.Text = "195 Lewis Road Suite" + [char]13 + "San Jose, CA 95111"
}

# 3. Remove the blank "NoSpacing" paragraph that immediately follows the
#    "Board of Directors" paragraph.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Board of Directors*") {
        $next = $paras.Item($i + 1)
        $next.Range.Delete()
        break
    }
}
